# "add antenna stick and add diagramm"
# Replace the 36-point (0..350 step 10) antenna-pattern table with a denser
# 72-point (0..355 step 5) table, re-point the radar chart's series at the
# new ranges, rebase the value axis minimum to 0, nudge the chart's anchor
# (as Excel itself does when the plotted data range changes), and leave the
# selection on B2 like the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New angle (A) / level (B) data, rows 2-73 --------------------------
$aVals = @(0,5,10,15,20,25,30,35,40,45,50,55,60,65,70,75,80,85,90,95,100,105,110,115,120,125,130,135,140,145,150,155,160,165,170,175,180,185,190,195,200,205,210,215,220,225,230,235,240,245,250,255,260,265,270,275,280,285,290,295,300,305,310,315,320,325,330,335,340,345,350,355)
$bVals = @(0,12.7,13,12.7,12,10,3,3,3,3,3,3,3,3,0,0,0,5,0,0,0,0,0,0,0,0,0,0,0,3,0,0,0,0,0,0,0,0,3,0,0,0,0,0,0,0,0,3,0,0,0,0,0,0,0,0,0,0,0,5,0,0,0,3,3,3,3,3,3,3,3,10)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

# Rows 38-73 are brand-new rows; give them the same centered style the
# original A2:B37 block used (xlCenter == -4108).
$ws.Range("A38:B73").HorizontalAlignment = -4108

# B2 loses its inherited "centered" cell style in the source edit.
$ws.Range("B2").ClearFormats()

# ---- Chart: repoint series ranges + rescale value axis -------------------
$co = $ws.ChartObjects().Item(1)
$ch = $co.Chart
$s = $ch.SeriesCollection().Item(1)
$s.Formula = '=SERIES(,Лист1!$A$2:$A$73,Лист1!$B$4:$B$73,1)'

$valAxis = $ch.Axes(2)
$valAxis.MinimumScale = 0

# ---- Chart anchor nudges to match the re-laid-out diagram ----------------
$co.Left = 235.25
$co.Top = 26.25
$co.Width = 681.8125
$co.Height = 402

# ---- Misc view state -------------------------------------------------------
$ws.Range("B2").Select()
$ws.PageSetup.Orientation = 1
